$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the updated Price cells are digit strings that parse as
# plain numbers (e.g. "0.613"). Excel auto-converts Range.Value to a
# Double in that case, same as typing it into a General-formatted cell,
# which would turn the column back into numbers instead of text. Force
# those specific cells to Text first so the literal string sticks, then
# drop the formatting back to the default style once the value is set
# (all these share one Text style object, so this does not fan out).
$textCells = @("D5", "D6", "D7", "D9", "D13", "D14", "D19", "D22", "D26", "D27", "D29", "D30", "D32", "D33", "D36", "D40", "D42", "D44", "D45", "D46", "D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "38.274.94"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").Value = "2.094.01"
$ws.Range("E3").Value = "  +3.40%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "228.41"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "0.613"
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("D7").Value = "60.94"
$ws.Range("E7").Value = "  +2.02%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "0.379"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("E10").Value = "  +2.35%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "2.403.96"
$ws.Range("E12").Value = "  +3.30%  "
$ws.Range("D13").Value = "14.75"
$ws.Range("E13").Value = "  +2.66%  "
$ws.Range("D14").Value = "22.31"
$ws.Range("E14").Value = "  +6.41%  "
$ws.Range("E15").Value = "  +5.79%  "
$ws.Range("E16").Value = "  +2.60%  "
$ws.Range("D17").Value = "2.092.54"
$ws.Range("E17").Value = "  +3.26%  "
$ws.Range("D18").Value = "38.220.31"
$ws.Range("E18").Value = "  +1.61%  "
$ws.Range("D19").Value = "70.36"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").Value = "224.92"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +2.71%  "
$ws.Range("E25").Value = "  +3.21%  "
$ws.Range("D26").Value = "169.84"
$ws.Range("E26").Value = "  +1.49%  "
$ws.Range("D27").Value = "9.40"
$ws.Range("E27").Value = "  +2.03%  "
$ws.Range("E28").Value = "  +1.54%  "
$ws.Range("D29").Value = "18.99"
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("D30").Value = "1.37"
$ws.Range("E30").Value = "  +9.12%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "2.35"
$ws.Range("E32").Value = "  +6.03%  "
$ws.Range("D33").Value = "4.75"
$ws.Range("E33").Value = "  +6.64%  "
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("E35").Value = "  +0.39%  "
$ws.Range("D36").Value = "6.43"
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("E37").Value = "  +4.54%  "
$ws.Range("E38").Value = "  +5.02%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").Value = "18.20"
$ws.Range("E40").Value = "  +2.40%  "
$ws.Range("D41").Value = "1.543.01"
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D42").Value = "99.80"
$ws.Range("E42").Value = "  +4.83%  "
$ws.Range("E43").Value = "  +1.87%  "
$ws.Range("D44").Value = "2.83"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("D45").Value = "0.0909"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "4.14"
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("E47").Value = "  +1.22%  "
$ws.Range("D48").Value = "7.49"
$ws.Range("E48").Value = "  +5.96%  "
$ws.Range("E49").Value = "  +3.02%  "
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").Value = "2.290.10"
$ws.Range("E51").Value = "  +3.28%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
